$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as plain text (many use a European-style
# thousands separator such as "31.053.61", which is not a valid Excel number).
# Force NumberFormat to Text on the cells we are about to rewrite so Excel does
# not silently reinterpret the literal strings as numbers.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D20").NumberFormat = "@"
$ws.Range("D22:D26").NumberFormat = "@"
$ws.Range("D28:D30").NumberFormat = "@"
$ws.Range("D32:D34").NumberFormat = "@"
$ws.Range("D36:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.053.61'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '1.957.51'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '246.37'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4890'
$ws.Range("E7").Value = '  +1.50%  '
$ws.Range("D8").Value = '0.2974'
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = '0.06873'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").Value = '19.28'
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").Value = '108.16'
$ws.Range("E11").Value = '  -3.44%  '
$ws.Range("D12").Value = '1.959.43'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '0.07757'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '5.471'
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = '0.7097'
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("D16").Value = '283.68'
$ws.Range("E16").Value = '  -3.75%  '
$ws.Range("D17").Value = '31.076.04'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '13.30'
$ws.Range("E18").Value = '  +0.66%  '
$ws.Range("D19").Value = '0.000007777'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").Value = '2.203.65'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '5.537'
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '6.508'
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").Value = '9.838'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("D26").Value = '169.73'
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = '2.232'
$ws.Range("E28").Value = '  +2.41%  '
$ws.Range("D29").Value = '0.1058'
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("D30").Value = '1.434'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '4.610'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").Value = '4.466'
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = '0.04994'
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").Value = '1.183'
$ws.Range("E36").Value = '  +2.27%  '
$ws.Range("D37").Value = '2.728'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = '0.02033'
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").Value = '2.703'
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").Value = '2.177'
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("D41").Value = '6.500'
$ws.Range("E41").Value = '  +9.93%  '
$ws.Range("D42").Value = '75.71'
$ws.Range("E42").Value = '  +8.82%  '
$ws.Range("D43").Value = '0.4534'
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("D44").Value = '0.8874'
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").Value = '109.46'
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").Value = '8.135'
$ws.Range("E46").Value = '  +10.72%  '
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").Value = '991.97'
$ws.Range("E48").Value = '  +9.86%  '
$ws.Range("D49").Value = '9.444'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").Value = '0.1266'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = '0.2591'
$ws.Range("E51").Value = '  +2.99%  '
